$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new "Dataset" column before the existing G column (old G/H shift to H/I)
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").Value = "Dataset"
$ws.Range("G2:G7").Value = "clean"

# New run results rows
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 32
$ws.Range("D8").Value = 0.4492
$ws.Range("E8").Value = 0.3988
$ws.Range("F8").Value = 2.6706
$ws.Range("G8").Value = "cleaner"
$ws.Range("H8").Value = "the role of the earth in the solar system"
$ws.Range("I8").Value = "the size distribution of the neptune trojans and the missing intermediate sized planetesimals"

$ws.Range("A9").Value = 48
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 128
$ws.Range("D9").Value = 0.5191
$ws.Range("E9").Value = 0.4609
$ws.Range("F9").Value = 1.7501
$ws.Range("G9").Value = "cleaner"
$ws.Range("H9").Value = "the precis imag of the earth of the earth"
$ws.Range("I9").Value = "the size distribution of the neptune trojans and the missing intermediate sized planetesimals"

$ws.Range("A10").Formula = "=48*4"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 128
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.4512
$ws.Range("F10").Value = 0.0383
$ws.Range("G10").Value = "cleaner"
$ws.Range("H10").Value = "first space imag from a dslr object from the new planet <EOS>"
$ws.Range("I10").Value = "the size distribution of the neptune trojans and the missing intermediate sized planetesimals"

# Update selection to match final state
[void]$ws.Range("C11").Select()

# Touch page setup so a pageSetup element is emitted (portrait orientation)
$ws.PageSetup.Orientation = 1
